$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUJIY")

# Row 8
$ws.Range("D8").Value = 21997600
$ws.Range("E8").Value = 20992400
$ws.Range("F8").Value = 22241900
$ws.Range("G8").Value = 22533100
$ws.Range("H8").Value = 22057200
$ws.Range("I8").Value = 20020900
$ws.Range("J8").Value = 19845400

# Row 9
$ws.Range("D9").Value = 13228700
$ws.Range("E9").Value = 12579300
$ws.Range("F9").Value = 13485500
$ws.Range("G9").Value = 13841100
$ws.Range("H9").Value = 13731000
$ws.Range("I9").Value = 12332300
$ws.Range("J9").Value = 12001200

# Row 10
$ws.Range("D10").Value = 8768900
$ws.Range("E10").Value = 8413000
$ws.Range("F10").Value = 8756400
$ws.Range("G10").Value = 8692000
$ws.Range("H10").Value = 8326200
$ws.Range("I10").Value = 7688500
$ws.Range("J10").Value = 7844200

# Row 12
$ws.Range("D12").Value = 1503600
$ws.Range("E12").Value = 1448500
$ws.Range("F12").Value = 1473800
$ws.Range("G12").Value = 1456700
$ws.Range("H12").Value = 1485700
$ws.Range("I12").Value = 1520100
$ws.Range("J12").Value = 1567300

# Row 14
$ws.Range("J14").Value = 175400

# Row 17
$ws.Range("D17").Value = 20816300
$ws.Range("E17").Value = 19434900
$ws.Range("F17").Value = 20609000
$ws.Range("G17").Value = 20974700
$ws.Range("H17").Value = 20784300
$ws.Range("I17").Value = 18989200
$ws.Range("J17").Value = 18999800

# Row 18
$ws.Range("D18").Value = 1181300
$ws.Range("E18").Value = 1557400
$ws.Range("F18").Value = 1632900
$ws.Range("G18").Value = 1558500
$ws.Range("H18").Value = 1272900
$ws.Range("I18").Value = 1031600
$ws.Range("J18").Value = 845700

# Row 20
$ws.Range("D20").Value = 648100
$ws.Range("E20").Value = 246700
$ws.Range("F20").Value = 54200
$ws.Range("G20").Value = 264600
$ws.Range("H20").Value = 185600
$ws.Range("I20").Value = 85300
$ws.Range("J20").Value = -8500

# Row 21
$ws.Range("D21").Value = 2984600
$ws.Range("E21").Value = 2868700
$ws.Range("F21").Value = 2817900
$ws.Range("G21").Value = 2935400
$ws.Range("H21").Value = 2738500
$ws.Range("I21").Value = 2397900
$ws.Range("J21").Value = 2174500

# Row 22
$ws.Range("D22").Value = 41300
$ws.Range("E22").Value = 43300
$ws.Range("F22").Value = 39600
$ws.Range("G22").Value = 41300
$ws.Range("H22").Value = 37800
$ws.Range("I22").Value = 39400
$ws.Range("J22").Value = 30900

# Row 23
$ws.Range("D23").Value = 1788200
$ws.Range("E23").Value = 1760800
$ws.Range("F23").Value = 1647500
$ws.Range("G23").Value = 1781800
$ws.Range("H23").Value = 1420700
$ws.Range("I23").Value = 1077400
$ws.Range("J23").Value = 806300

# Row 24
$ws.Range("D24").Value = 492100
$ws.Range("E24").Value = 398900
$ws.Range("F24").Value = 437900
$ws.Range("G24").Value = 528700
$ws.Range("H24").Value = 494500
$ws.Range("I24").Value = 404400
$ws.Range("J24").Value = 271500

# Row 26
$ws.Range("D26").Value = 1296100
$ws.Range("E26").Value = 1361900
$ws.Range("F26").Value = 1209500
$ws.Range("G26").Value = 1253100
$ws.Range("H26").Value = 926200
$ws.Range("I26").Value = 673000
$ws.Range("J26").Value = 534800

# Row 27
$ws.Range("D27").Value = 1271900
$ws.Range("E27").Value = 1188800
$ws.Range("F27").Value = 1052300
$ws.Range("G27").Value = 1071700
$ws.Range("H27").Value = 732200
$ws.Range("I27").Value = 490600
$ws.Range("J27").Value = 395600

# Row 32
$ws.Range("D32").Value = -648100
$ws.Range("E32").Value = -246700
$ws.Range("F32").Value = -54200
$ws.Range("G32").Value = -264600
$ws.Range("H32").Value = -185600
$ws.Range("I32").Value = -85300
$ws.Range("J32").Value = 8500

# Row 33
$ws.Range("D33").Value = 1271900
$ws.Range("E33").Value = 1188800
$ws.Range("F33").Value = 1052300
$ws.Range("G33").Value = 1071700
$ws.Range("H33").Value = 732200
$ws.Range("I33").Value = 490600
$ws.Range("J33").Value = 395600

# Row 35
$ws.Range("D35").Value = 1271900
$ws.Range("E35").Value = 1188800
$ws.Range("F35").Value = 1052300
$ws.Range("G35").Value = 1071700
$ws.Range("H35").Value = 732200
$ws.Range("I35").Value = 490600
$ws.Range("J35").Value = 395600

# Row 41
$ws.Range("D41").Value = 6944900
$ws.Range("E41").Value = 7918700
$ws.Range("F41").Value = 5432100
$ws.Range("G41").Value = 6571100
$ws.Range("H41").Value = 5465300
$ws.Range("I41").Value = 4026400
$ws.Range("J41").Value = 2125300

# Row 42
$ws.Range("E42").Value = 54600
$ws.Range("F42").Value = 253200
$ws.Range("G42").Value = 172100
$ws.Range("H42").Value = 150400
$ws.Range("I42").Value = 33000
$ws.Range("J42").Value = 111800

# Row 43
$ws.Range("D43").Value = 5596200
$ws.Range("E43").Value = 5747500
$ws.Range("F43").Value = 5817900
$ws.Range("G43").Value = 6156400
$ws.Range("H43").Value = 5756800
$ws.Range("I43").Value = 5323400
$ws.Range("J43").Value = 5030600

# Row 44
$ws.Range("D44").Value = 3266300
$ws.Range("E44").Value = 3066700
$ws.Range("F44").Value = 3156500
$ws.Range("G44").Value = 3367500
$ws.Range("H44").Value = 3288200
$ws.Range("I44").Value = 3615400
$ws.Range("J44").Value = 3416700

# Row 45
$ws.Range("D45").Value = 1013800
$ws.Range("E45").Value = 1607100
$ws.Range("F45").Value = 1307200
$ws.Range("G45").Value = 1299700
$ws.Range("H45").Value = 1265000
$ws.Range("I45").Value = 1115900
$ws.Range("J45").Value = 1266400

# Row 46
$ws.Range("D46").Value = 16821700
$ws.Range("E46").Value = 18394600
$ws.Range("F46").Value = 15966900
$ws.Range("G46").Value = 17566700
$ws.Range("H46").Value = 15925700
$ws.Range("I46").Value = 14114000
$ws.Range("J46").Value = 11950800

# Row 47
$ws.Range("D47").Value = 2162000
$ws.Range("E47").Value = 2707800
$ws.Range("F47").Value = 2889900
$ws.Range("G47").Value = 3443500
$ws.Range("H47").Value = 3005400
$ws.Range("I47").Value = 2911400
$ws.Range("J47").Value = 2558900

# Row 48
$ws.Range("D48").Value = 4859400
$ws.Range("E48").Value = 4706400
$ws.Range("F48").Value = 4814700
$ws.Range("G48").Value = 4767500
$ws.Range("H48").Value = 4793300
$ws.Range("I48").Value = 4937000
$ws.Range("J48").Value = 5007400

# Row 49
$ws.Range("D49").Value = 6497100
$ws.Range("E49").Value = 5209600
$ws.Range("F49").Value = 5361800
$ws.Range("G49").Value = 5290500
$ws.Range("H49").Value = 4569600
$ws.Range("I49").Value = 4542500
$ws.Range("J49").Value = 3954500

# Row 52
$ws.Range("D52").Value = 1236100
$ws.Range("E52").Value = 921700
$ws.Range("F52").Value = 906900
$ws.Range("G52").Value = 1083200
$ws.Range("H52").Value = 877700
$ws.Range("I52").Value = 1153800
$ws.Range("J52").Value = 1295000

# Row 54
$ws.Range("D54").Value = 31576200
$ws.Range("E54").Value = 31940000
$ws.Range("F54").Value = 29940200
$ws.Range("G54").Value = 32151400
$ws.Range("H54").Value = 29171800
$ws.Range("I54").Value = 27658700
$ws.Range("J54").Value = 24766600

# Row 57
$ws.Range("D57").Value = 2055500
$ws.Range("E57").Value = 2151100
$ws.Range("F57").Value = 2112900
$ws.Range("G57").Value = 2280300
$ws.Range("H57").Value = 2245900
$ws.Range("I57").Value = 2081800
$ws.Range("J57").Value = 2336000

# Row 58
$ws.Range("D58").Value = 376800
$ws.Range("E58").Value = 1121000
$ws.Range("F58").Value = 500000
$ws.Range("G58").Value = 331300
$ws.Range("H58").Value = 404400
$ws.Range("I58").Value = 367700
$ws.Range("J58").Value = 1614000

# Row 59
$ws.Range("D59").Value = 3572300
$ws.Range("E59").Value = 2876800
$ws.Range("F59").Value = 2872700
$ws.Range("G59").Value = 3013900
$ws.Range("H59").Value = 2741400
$ws.Range("I59").Value = 2489600
$ws.Range("J59").Value = 2309100

# Row 60
$ws.Range("D60").Value = 6004500
$ws.Range("E60").Value = 6148800
$ws.Range("F60").Value = 5485500
$ws.Range("G60").Value = 5625500
$ws.Range("H60").Value = 5391700
$ws.Range("I60").Value = 4939100
$ws.Range("J60").Value = 6259000

# Row 61
$ws.Range("D61").Value = 3729000
$ws.Range("E61").Value = 3931000
$ws.Range("F61").Value = 2805900
$ws.Range("G61").Value = 2829900
$ws.Range("H61").Value = 2847300
$ws.Range("I61").Value = 2871000
$ws.Range("J61").Value = 183800

# Row 62
$ws.Range("D62").Value = 1062300
$ws.Range("E62").Value = 1357000
$ws.Range("F62").Value = 1471500
$ws.Range("G62").Value = 1390600
$ws.Range("H62").Value = 1060900
$ws.Range("I62").Value = 1544500
$ws.Range("J62").Value = 1541100

# Row 66
$ws.Range("D66").Value = 12780800
$ws.Range("E66").Value = 13466300
$ws.Range("F66").Value = 11726200
$ws.Range("G66").Value = 11967600
$ws.Range("H66").Value = 10905200
$ws.Range("I66").Value = 10764200
$ws.Range("J66").Value = 9201800

# Row 72
$ws.Range("D72").Value = 21549500
$ws.Range("E72").Value = 20571700
$ws.Range("F72").Value = 19661700
$ws.Range("G72").Value = 19219700
$ws.Range("H72").Value = 18409500
$ws.Range("I72").Value = 17895200
$ws.Range("J72").Value = 17578800

# Row 76
$ws.Range("D76").Value = 18795400
$ws.Range("E76").Value = 18473800
$ws.Range("F76").Value = 18214000
$ws.Range("G76").Value = 20183700
$ws.Range("H76").Value = 18266600
$ws.Range("I76").Value = 16894600
$ws.Range("J76").Value = 15564800

# Row 81
$ws.Range("D81").Value = 1271900
$ws.Range("E81").Value = 1188800
$ws.Range("F81").Value = 1052300
$ws.Range("G81").Value = 1071700
$ws.Range("H81").Value = 732200
$ws.Range("I81").Value = 490600
$ws.Range("J81").Value = 395600

# Row 83
$ws.Range("D83").Value = 1153900
$ws.Range("E83").Value = 1063400
$ws.Range("F83").Value = 1129600
$ws.Range("G83").Value = 1111100
$ws.Range("H83").Value = 1278600
$ws.Range("I83").Value = 1279600
$ws.Range("J83").Value = 1335900

# Row 89
$ws.Range("D89").Value = 2360800
$ws.Range("E89").Value = 2609100
$ws.Range("F89").Value = 2020300
$ws.Range("G89").Value = 2384100
$ws.Range("H89").Value = 2644500
$ws.Range("I89").Value = 1803000
$ws.Range("J89").Value = 1221600

# Row 91
$ws.Range("D91").Value = -569900
$ws.Range("E91").Value = -674800
$ws.Range("F91").Value = -576800
$ws.Range("G91").Value = -514800
$ws.Range("H91").Value = -635400
$ws.Range("I91").Value = -812500
$ws.Range("J91").Value = -910900

# Row 94
$ws.Range("D94").Value = -1010500
$ws.Range("E94").Value = -1052600
$ws.Range("F94").Value = -1422200
$ws.Range("G94").Value = -1089400
$ws.Range("H94").Value = -1134300
$ws.Range("I94").Value = -1274000
$ws.Range("J94").Value = -1680300

# Row 96
$ws.Range("D96").Value = -286700
$ws.Range("E96").Value = -272700
$ws.Range("F96").Value = -289000
$ws.Range("G96").Value = -239700
$ws.Range("H96").Value = -174200
$ws.Range("I96").Value = -163300
$ws.Range("J96").Value = -141500

# Row 100
$ws.Range("D100").Value = -2341000
$ws.Range("E100").Value = 1006100
$ws.Range("F100").Value = -1551900
$ws.Range("G100").Value = -412200
$ws.Range("H100").Value = -226800
$ws.Range("I100").Value = 1159700
$ws.Range("J100").Value = -220600

# Row 101
$ws.Range("D101").Value = 17000
$ws.Range("E101").Value = -76000
$ws.Range("F101").Value = -185200
$ws.Range("G101").Value = 223200
$ws.Range("H101").Value = 155600
$ws.Range("I101").Value = 212300
$ws.Range("J101").Value = -25500

# Row 102
$ws.Range("D102").Value = -973700
$ws.Range("E102").Value = 2486600
$ws.Range("F102").Value = -1139000
$ws.Range("G102").Value = 1105700
$ws.Range("H102").Value = 1439000
$ws.Range("I102").Value = 1901000
$ws.Range("J102").Value = -704800
